$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: square of column C (deviation from mean), mirroring the
# existing pattern used for column B (square of column A).
$ws.Range("D2").Formula = "=C2^2"
$ws.Range("D3:D11").Formula = "=C3^2"
$ws.Range("D12").Formula = "=AVERAGE(D2:D11)"

# Match the author's final selection being on the newly added D12 cell.
[void]$ws.Range("D12").Select()
